# Auto-generated Excel COM-interop script
# Applies numeric value updates to the Malboro Profits workbook sheets
# (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) to match the target diff.

$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H69").Value = 8249.5
$ws.Range("I69").Value = 7666
$ws.Range("K69").Value = 22998
$ws.Range("M69").Value = -22124
$ws.Range("H72").Value = 8249.5
$ws.Range("I72").Value = 7666
$ws.Range("K72").Value = 68994
$ws.Range("M72").Value = -64626
$ws.Range("H113").Value = 12542.8
$ws.Range("I113").Value = 20466.8
$ws.Range("J113").Value = 7260.1333
$ws.Range("K113").Value = 20466.8
$ws.Range("L113").Value = 7260.1333
$ws.Range("M113").Value = -17212.8
$ws.Range("N113").Value = -13768.1333
$ws.Range("H116").Value = 4846.7646
$ws.Range("I116").Value = 4933.1665
$ws.Range("J116").Value = 4639.4
$ws.Range("K116").Value = 4933.1665
$ws.Range("L116").Value = 4639.4
$ws.Range("M116").Value = -1491.1665
$ws.Range("N116").Value = -11523.4
$ws.Range("H132").Value = 6431.327
$ws.Range("I132").Value = 5006.9785
$ws.Range("J132").Value = 14799.375
$ws.Range("K132").Value = 15020.9355
$ws.Range("L132").Value = 44398.125
$ws.Range("M132").Value = -12490.9355
$ws.Range("N132").Value = -49458.125
$ws.Range("H133").Value = 64089
$ws.Range("J133").Value = 64089
$ws.Range("L133").Value = 64089
$ws.Range("N133").Value = -74209
$ws.Range("H138").Value = 2988.6487
$ws.Range("I138").Value = 3551.5
$ws.Range("K138").Value = 10654.5
$ws.Range("M138").Value = -5514.5
$ws.Range("H141").Value = 1037.1364
$ws.Range("I141").Value = 493.41177
$ws.Range("J141").Value = 2885.8
$ws.Range("K141").Value = 1480.23531
$ws.Range("L141").Value = 8657.400000000001
$ws.Range("M141").Value = 3699.76469
$ws.Range("N141").Value = -19017.4

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6513.7646
$ws.Range("I32").Value = 5802.8184
$ws.Range("K32").Value = 5802.8184
$ws.Range("M32").Value = -5515.8184
$ws.Range("H34").Value = 40000
$ws.Range("J34").Value = 40000
$ws.Range("L34").Value = 40000
$ws.Range("N34").Value = -40542
$ws.Range("H61").Value = 17125.705
$ws.Range("I61").Value = 5431.1
$ws.Range("J61").Value = 33832.285
$ws.Range("K61").Value = 5431.1
$ws.Range("L61").Value = 33832.285
$ws.Range("M61").Value = -5219.1
$ws.Range("N61").Value = -34256.285
$ws.Range("H74").Value = 20909.773
$ws.Range("I74").Value = 3017.5
$ws.Range("J74").Value = 27619.375
$ws.Range("K74").Value = 3017.5
$ws.Range("L74").Value = 27619.375
$ws.Range("M74").Value = -2143.5
$ws.Range("N74").Value = -29367.375
$ws.Range("H77").Value = 20909.773
$ws.Range("I77").Value = 3017.5
$ws.Range("J77").Value = 27619.375
$ws.Range("K77").Value = 15087.5
$ws.Range("L77").Value = 138096.875
$ws.Range("M77").Value = -10719.5
$ws.Range("N77").Value = -146832.875
$ws.Range("H93").Value = 18187.5
$ws.Range("I93").Value = 11428.571
$ws.Range("J93").Value = 65500
$ws.Range("K93").Value = 11428.571
$ws.Range("L93").Value = 65500
$ws.Range("M93").Value = -8932.571
$ws.Range("N93").Value = -70492
$ws.Range("H136").Value = 17125.705
$ws.Range("I136").Value = 5431.1
$ws.Range("J136").Value = 33832.285
$ws.Range("K136").Value = 16293.3
$ws.Range("L136").Value = 101496.855
$ws.Range("M136").Value = -13743.3
$ws.Range("N136").Value = -106596.855
$ws.Range("H140").Value = 89400
$ws.Range("J140").Value = 89400
$ws.Range("L140").Value = 89400
$ws.Range("N140").Value = -99760

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H38").Value = 0
$ws.Range("J38").Value = 0
$ws.Range("L38").Value = 0
$ws.Range("N38").ClearContents()
$ws.Range("H82").Value = 1064.6666
$ws.Range("I82").Value = 1064.6666
$ws.Range("J82").Value = 0
$ws.Range("K82").Value = 1064.6666
$ws.Range("L82").Value = 0
$ws.Range("M82").Value = -681.6666
$ws.Range("N82").ClearContents()
$ws.Range("H85").Value = 1064.6666
$ws.Range("I85").Value = 1064.6666
$ws.Range("J85").Value = 0
$ws.Range("K85").Value = 1064.6666
$ws.Range("L85").Value = 0
$ws.Range("M85").Value = 261.3334
$ws.Range("N85").ClearContents()

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 36097.734
$ws.Range("I31").Value = 51666
$ws.Range("J31").Value = 33702.617
$ws.Range("K31").Value = 51666
$ws.Range("L31").Value = 33702.617
$ws.Range("M31").Value = -51371
$ws.Range("N31").Value = -34292.617
$ws.Range("H34").Value = 36097.734
$ws.Range("I34").Value = 51666
$ws.Range("J34").Value = 33702.617
$ws.Range("K34").Value = 51666
$ws.Range("L34").Value = 33702.617
$ws.Range("M34").Value = -51464
$ws.Range("N34").Value = -34106.617
$ws.Range("H58").Value = 16408.793
$ws.Range("J58").Value = 16962.625
$ws.Range("L58").Value = 16962.625
$ws.Range("N58").Value = -17368.625
$ws.Range("H136").Value = 16408.793
$ws.Range("J136").Value = 16962.625
$ws.Range("L136").Value = 50887.875
$ws.Range("N136").Value = -55987.875
$ws.Range("H140").Value = 74605.16
$ws.Range("J140").Value = 74605.16
$ws.Range("L140").Value = 74605.16
$ws.Range("N140").Value = -84965.16

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H76").Value = 4553
$ws.Range("I76").Value = 4553
$ws.Range("K76").Value = 13659
$ws.Range("M76").Value = -13276
$ws.Range("H79").Value = 4553
$ws.Range("I79").Value = 4553
$ws.Range("K79").Value = 13659
$ws.Range("M79").Value = -12333
$ws.Range("H131").Value = 1437.4849
$ws.Range("I131").Value = 922.2857
$ws.Range("J131").Value = 1476.6848
$ws.Range("K131").Value = 2766.8571
$ws.Range("L131").Value = 4430.0544
$ws.Range("M131").Value = 2273.1429
$ws.Range("N131").Value = -14510.0544

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H39").Value = 36995.383
$ws.Range("J39").Value = 36995.383
$ws.Range("L39").Value = 36995.383
$ws.Range("N39").Value = -38059.383
$ws.Range("H52").Value = 29161.445
$ws.Range("J52").Value = 29161.445
$ws.Range("L52").Value = 29161.445
$ws.Range("N52").Value = -29679.445
$ws.Range("H93").Value = 36909.656
$ws.Range("J93").Value = 36909.656
$ws.Range("L93").Value = 36909.656
$ws.Range("N93").Value = -40653.656
$ws.Range("H113").Value = 4795.909
$ws.Range("I113").Value = 3393
$ws.Range("K113").Value = 3393
$ws.Range("M113").Value = -1223
$ws.Range("H132").Value = 21392
$ws.Range("I132").Value = 12357.143
$ws.Range("K132").Value = 37071.429
$ws.Range("M132").Value = -34541.429

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H29").Value = 2000
$ws.Range("J29").Value = 2000
$ws.Range("L29").Value = 2000
$ws.Range("N29").Value = -2590
$ws.Range("H33").Value = 0
$ws.Range("J33").Value = 0
$ws.Range("L33").Value = 0
$ws.Range("N33").ClearContents()
$ws.Range("H35").Value = 2685.3333
$ws.Range("I35").Value = 644.25
$ws.Range("K35").Value = 644.25
$ws.Range("M35").Value = -308.25
$ws.Range("H97").Value = 30000
$ws.Range("J97").Value = 30000
$ws.Range("L97").Value = 30000
$ws.Range("N97").Value = -31982
$ws.Range("H132").Value = 2362836.5
$ws.Range("I132").Value = 3991.8
$ws.Range("K132").Value = 11975.4
$ws.Range("M132").Value = -9445.400000000001
$ws.Range("H136").Value = 12730.286
$ws.Range("J136").Value = 12300.0625
$ws.Range("L136").Value = 36900.1875
$ws.Range("N136").Value = -42000.1875
$ws.Range("H138").Value = 250000
$ws.Range("J138").Value = 250000
$ws.Range("L138").Value = 250000
$ws.Range("N138").Value = -260280

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H11").Value = 5750
$ws.Range("J11").Value = 5750
$ws.Range("L11").Value = 5750
$ws.Range("N11").Value = -6034
$ws.Range("H51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("N51").ClearContents()
$ws.Range("H52").Value = 25344.8
$ws.Range("H122").Value = 3737.7021
$ws.Range("I122").Value = 1725.2858
$ws.Range("J122").Value = 6703.3687
$ws.Range("K122").Value = 5175.857400000001
$ws.Range("L122").Value = 20110.1061
$ws.Range("M122").Value = -2725.857400000001
$ws.Range("N122").Value = -25010.1061
$ws.Range("H124").Value = 47500
$ws.Range("J124").Value = 47500
$ws.Range("L124").Value = 47500
$ws.Range("N124").Value = -57320
$ws.Range("H132").Value = 8311.5
$ws.Range("I132").Value = 936.1111
$ws.Range("J132").Value = 41500.75
$ws.Range("K132").Value = 2808.3333
$ws.Range("L132").Value = 124502.25
$ws.Range("M132").Value = -278.3332999999998
$ws.Range("N132").Value = -129562.25
$ws.Range("H136").Value = 13167.333
$ws.Range("I136").Value = 2904.6365
$ws.Range("J136").Value = 24456.3
$ws.Range("K136").Value = 8713.9095
$ws.Range("L136").Value = 73368.89999999999
$ws.Range("M136").Value = -6163.9095
$ws.Range("N136").Value = -78468.89999999999
